# Task Management.xlsx - "moving data under resources"
#
# Real content edits made by the author (everything else in the diff is
# cosmetic noise from saving with a different Excel build/locale):
#   1. B4: the "For LA yelp data ..." task text was reworded to refer to
#      "Las Vegas" data instead of "LA" data.
#   2. E7: the "Done" progress marker was cleared out entirely.
#   3. D12: the "No of people required" count was bumped from 1 to 2.
#   4. The user's last selection before saving was cell B8 (previously E16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "For Las Vegas yelp data, Find how much missing data exist - which restaurants have all data missing, some data missing "
$ws.Range("D12").Value = 2
$ws.Range("E7").Clear() | Out-Null

$ws.Range("B8").Select() | Out-Null
